# Update Maryland overview workbook: convert key count figures to
# text-formatted numbers (with thousands separators where applicable)
# and append a "Total" row to the County sheet.

$wb = $excel.ActiveWorkbook

# Writes a vertical run of text values down a single column, starting
# at ($startRow, $col).
function Set-TextValuesDown {
    param($ws, $rangeAddr, $startRow, $col, $vals)
    $rng = $ws.Range($rangeAddr)
    $rng.NumberFormat = "@"
    for ($i = 0; $i -lt $vals.Count; $i++) {
        $ws.Cells.Item($startRow + $i, $col).Value = $vals[$i]
    }
    $rng.ClearFormats()
}

# Writes a horizontal run of text values across a single row, starting
# at ($row, $startCol).
function Set-TextValuesAcross {
    param($ws, $rangeAddr, $row, $startCol, $vals)
    $rng = $ws.Range($rangeAddr)
    $rng.NumberFormat = "@"
    for ($i = 0; $i -lt $vals.Count; $i++) {
        $ws.Cells.Item($row, $startCol + $i).Value = $vals[$i]
    }
    $rng.ClearFormats()
}

# ---------------------------------------------------------------------
# Sheet "Overall": A2 total filer count -> text "2,487"
# ---------------------------------------------------------------------
$wsOverall = $wb.Worksheets.Item("Overall")
Set-TextValuesAcross $wsOverall "A2" 2 1 @("2,487")

# ---------------------------------------------------------------------
# Sheet "County": per-county filer counts -> text; add Total row 26
# ---------------------------------------------------------------------
$wsCounty = $wb.Worksheets.Item("County")
$countyCounts = @("42","195","236","447","29","20","52","34","53","20","92","20","81","126","21","513","196","26","15","40","45","96","61","27")
Set-TextValuesDown $wsCounty "B2:B25" 2 2 $countyCounts

Set-TextValuesAcross $wsCounty "A26:F26" 26 1 @("Total", "2,487", '$6,127,216,573', "8.58%", "-15.15%", "69.08%")

# ---------------------------------------------------------------------
# Sheet "Congressional District": per-district filer counts -> text;
# Total row (B11) -> text "2,487"
# ---------------------------------------------------------------------
$wsCD = $wb.Worksheets.Item("Congressional District")
$cdCounts = @("362","273","297","156","191","304","448","452","4","2,487")
Set-TextValuesDown $wsCD "B2:B11" 2 2 $cdCounts

# ---------------------------------------------------------------------
# Sheet "Size": per-size-bucket filer counts -> text;
# Total row (B8) -> text "2,487"
# ---------------------------------------------------------------------
$wsSize = $wb.Worksheets.Item("Size")
$sizeCounts = @("776","679","352","189","357","134","2,487")
Set-TextValuesDown $wsSize "B2:B8" 2 2 $sizeCounts

# ---------------------------------------------------------------------
# Sheet "Subsector": per-subsector filer counts -> text;
# Total row (B14) -> text "2,487"
# ---------------------------------------------------------------------
$wsSubsector = $wb.Worksheets.Item("Subsector")
$subsectorCounts = @("227","233","108","251","15","816","29","2","221","35","531","19","2,487")
Set-TextValuesDown $wsSubsector "B2:B14" 2 2 $subsectorCounts
